$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume number (44 -> 45) and report week dates ---
$ws.Range("A8").Value = "Volume 30   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/6/2023  Through  11/12/2023"

# --- Weekly crime statistics table (rows 16-30) ---
# Cells that must display the report's "no data" placeholder text (shared
# strings "0" / "***.*") are populated by copying from an existing placeholder
# cell so the exact text style (s=14) and shared-string reuse is preserved.

$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 45
$ws.Range("J16").Value = 54
$ws.Range("K16").Value = -16.666666666666
$ws.Range("L16").Value = 28.571428571428
$ws.Range("M16").Value = -52.127659574468
$ws.Range("N16").Value = -83.018867924528
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = -41.666666666666
$ws.Range("I17").Value = 148
$ws.Range("J17").Value = 114
$ws.Range("K17").Value = 29.824561403508
$ws.Range("L17").Value = 37.037037037037
$ws.Range("M17").Value = 20.325203252032
$ws.Range("N17").Value = -46.181818181818
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 77
$ws.Range("J18").Value = 68
$ws.Range("K18").Value = 13.235294117647
$ws.Range("L18").Value = 18.461538461538
$ws.Range("M18").Value = -58.152173913043
$ws.Range("N18").Value = -93.513058129738
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 21.428571428571
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -14.583333333333
$ws.Range("I19").Value = 418
$ws.Range("J19").Value = 316
$ws.Range("K19").Value = 32.278481012658
$ws.Range("L19").Value = 49.820788530465
$ws.Range("M19").Value = 16.759776536312
$ws.Range("N19").Value = -42.185338865836
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -50
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = -70
$ws.Range("I20").Value = 98
$ws.Range("J20").Value = 121
$ws.Range("K20").Value = -19.008264462809
$ws.Range("L20").Value = 40
$ws.Range("M20").Value = -5.769230769230
$ws.Range("N20").Value = -95.967078189300
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = -30.337078651685
$ws.Range("I21").Value = 796
$ws.Range("J21").Value = 682
$ws.Range("K21").Value = 16.715542521994
$ws.Range("L21").Value = 40.884955752212
$ws.Range("M21").Value = -9.648127128263
$ws.Range("N21").Value = -83.755102040816
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("C23").Value = 1
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("D23").Value = 1
$ws.Range("E23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 2
$ws.Range("I23").Value = 30
$ws.Range("J23").Value = 17
$ws.Range("K23").Value = 76.470588235294
$ws.Range("L23").Value = -14.285714285714
$ws.Range("M23").Value = 42.857142857142
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 13.333333333333
$ws.Range("F24").Value = 90
$ws.Range("G24").Value = 96
$ws.Range("H24").Value = -6.25
$ws.Range("I24").Value = 946
$ws.Range("J24").Value = 747
$ws.Range("K24").Value = 26.639892904953
$ws.Range("L24").Value = 118.981481481481
$ws.Range("M24").Value = -36.680053547523
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 5
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = 9.090909090909
$ws.Range("I25").Value = 279
$ws.Range("J25").Value = 292
$ws.Range("K25").Value = -4.452054794520
$ws.Range("L25").Value = 27.981651376146
$ws.Range("M25").Value = -40.889830508474
$ws.Range("D27").Copy($ws.Range("C27"))
$ws.Range("F27").Value = 7
$ws.Range("D27").Copy($ws.Range("G27"))
$ws.Range("E27").Copy($ws.Range("H27"))
$ws.Range("D27").Copy($ws.Range("F30"))
